# The deck originally ships two embedded themes:
#   ppt/theme/theme1.xml -> "Integral"      (used by the slide master / slides)
#   ppt/theme/theme2.xml -> "Office Theme"  (used by the notes master)
#
# The target revision swaps which palette the slide master uses: the
# slide-facing theme (theme1.xml, reached here through
# ActivePresentation.SlideMaster.Theme.ThemeColorScheme) becomes the
# stock "Office Theme" palette (the one formerly only used by the notes
# master / theme2.xml), while theme2.xml keeps the "Integral" colors.
#
# theme1.xml and theme2.xml already share an identical font scheme and
# format scheme (fills/lines/effects) -- the only real difference
# between the two theme parts is their 12 color-scheme entries (plus
# the cosmetic theme/clrScheme "name" attributes, which PowerPoint's
# object model does not expose for editing). So reproducing the new
# color values on the slide master's theme is the part of this change
# that is achievable -- and sufficient -- through the COM object model.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# Office Theme color scheme (previously theme2.xml), applied in the
# canonical clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$cs.Colors(1).RGB  = 0         # dk1      000000
$cs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
